# The "Email" config sheet had an extra "name"/"Delia" row that is no
# longer used, and the email address needs to be updated to the new
# recipient. Remove that row (shifting subject/body up) and update the
# email value.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Email")

# Row 3 holds the "name" / "Delia" / "Name of the person" entry - delete it
# entirely so the rows below (subject, body) shift up.
$ws.Rows.Item(3).Delete() | Out-Null

# Update the email address value (row 2, column B).
$ws.Cells.Item(2, 2).Value = "delia.panca@fwfcompany.com"

# Reflect the new selection left behind after the edit.
$ws.Range("B6").Select() | Out-Null
